$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old row 2 (data row), shifting rows 3..21 up to become rows 2..20
$ws.Rows.Item(2).Delete()

# Append the new data rows (21..31) that were added in this edit
$newData = @(
  @(4.037018136492105, -2.376052915471196, -2.018577781695872),
  @(5.506457075331898, -2.59009868774606, -0.0266988300582125),
  @(2.450691733545613, 1.12102667965644, -1.359096467929008),
  @(-3.604418606434045, 0.8684271347926146, -0.2955643468111461),
  @(-3.740167403684082, -2.81262268075148, -0.7813068606610404),
  @(-0.1918180868463897, 4.700260833629038, -2.898207216586866),
  @(3.375462932494072, 2.835833433762467, -3.255329332305403),
  @(3.974689759097054, -8.058607795863464, -3.702626524619664),
  @(1.054417246753779, -4.126515001348001, -1.594708632497015),
  @(-4.560519280942842, 2.439459521793583, 8.084467855471887),
  @(-1.622863769531381, 3.112303316593201, 3.16303658485438)
)

$startRow = 21
for ($i = 0; $i -lt $newData.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newData[$i][0]
    $ws.Cells.Item($r, 2).Value = $newData[$i][1]
    $ws.Cells.Item($r, 3).Value = $newData[$i][2]
}
